$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Pais" sheet tracks per-country COVID-19 stats; this refresh updates the
# snapshot timestamp and the case counts (columns B:H) that changed between
# the 12:52 and 13:22 pulls. Column A (country) is untouched per row.

$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 13:22"

# Row 11: Iran
$ws.Range("B11").Value = 71686
$ws.Range("C11").Value = 1657
$ws.Range("D11").Value = 43894
$ws.Range("E11").Value = 23318
$ws.Range("F11").Value = 3930
$ws.Range("G11").Value = 117
$ws.Range("H11").Value = 4474

# Row 30: Noruega
$ws.Range("B30").Value = 6409
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 32
$ws.Range("E30").Value = 6257
$ws.Range("F30").Value = 67
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 120

# Row 36: Pakistan
$ws.Range("B36").Value = 5170
$ws.Range("C36").Value = 159
$ws.Range("D36").Value = 1026
$ws.Range("E36").Value = 4056
$ws.Range("F36").Value = 37
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 88

# Row 46: Finlandia
$ws.Range("B46").Value = 2979
$ws.Range("C46").Value = 251
$ws.Range("D46").Value = 275
$ws.Range("E46").Value = 2697
$ws.Range("F46").Value = 37
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 7

# Row 47: Ucrania
$ws.Range("B47").Value = 2974
$ws.Range("C47").Value = 69
$ws.Range("D47").Value = 300
$ws.Range("E47").Value = 2618
$ws.Range("F47").Value = 80
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 56

# Row 48: Republica Dominicana
$ws.Range("B48").Value = 2777
$ws.Range("C48").Value = 266
$ws.Range("D48").Value = 89
$ws.Range("E48").Value = 2605
$ws.Range("F48").Value = 45
$ws.Range("G48").Value = 10
$ws.Range("H48").Value = 83

# Row 49: Catar
$ws.Range("B49").Value = 2759
$ws.Range("C49").Value = 0
$ws.Range("D49").Value = 108
$ws.Range("E49").Value = 2516
$ws.Range("F49").Value = 147
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 135

# Row 76: Camerun
$ws.Range("B76").Value = 833
$ws.Range("C76").Value = 66
$ws.Range("D76").Value = 42
$ws.Range("E76").Value = 787
$ws.Range("F76").Value = 8
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 4

# Row 77: Uzbekistan
$ws.Range("B77").Value = 828
$ws.Range("C77").Value = 68
$ws.Range("D77").Value = 41
$ws.Range("E77").Value = 753
$ws.Range("F77").Value = 15
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 34

# Row 78: Republica de Macedonia
$ws.Range("B78").Value = 820
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 98
$ws.Range("E78").Value = 710
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 12

# Row 114: Isla de Man
$ws.Range("B114").Value = 234
$ws.Range("C114").Value = 11
$ws.Range("D114").Value = 16
$ws.Range("E114").Value = 198
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 20

# Row 115: Consejo Danes para los Refugiados
$ws.Range("B115").Value = 226
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 112
$ws.Range("E115").Value = 112
$ws.Range("F115").Value = 11
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 2

# Row 167: Mozambique
$ws.Range("B167").Value = 21
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 3
$ws.Range("E167").Value = 18
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

# Row 168: Maldivas
$ws.Range("B168").Value = 20
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 2
$ws.Range("E168").Value = 18
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

# Row 169: Laos
$ws.Range("B169").Value = 20
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 13
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 0

# Row 170: Sudan
$ws.Range("B170").Value = 19
$ws.Range("C170").Value = 1
$ws.Range("D170").Value = 0
$ws.Range("E170").Value = 19
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 0

# Row 171: Angola
$ws.Range("B171").Value = 19
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 2
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 2

# Row 172: Nueva Caledonia
$ws.Range("B172").Value = 19
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 4
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 2

# Row 173: Guinea Ecuatorial
$ws.Range("B173").Value = 18
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 1
$ws.Range("E173").Value = 17
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
